# Updated run for publication
# Update computed frequency values in the FrequencyTables sheet (rows 2-5, cols B-X)
# to the latest values from the re-run analysis.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0754716981132075
$ws.Range("C2").Value = 0.00217706821480406
$ws.Range("E2").Value = 0.854862119013062
$ws.Range("F2").Value = 0.00145137880986938
$ws.Range("G2").Value = 0.000725689404934688
$ws.Range("H2").Value = 0.873730043541364
$ws.Range("J2").Value = 0.0362844702467344
$ws.Range("K2").Value = 0.000725689404934688
$ws.Range("L2").Value = 0.0297532656023222
$ws.Range("M2").Value = 0.831640058055152
$ws.Range("N2").Value = 0.00725689404934688
$ws.Range("O2").Value = 0.000725689404934688
$ws.Range("P2").Value = 0.991291727140784
$ws.Range("R2").Value = 0.421625544267054
$ws.Range("S2").Value = 0.0224963715529753
$ws.Range("T2").Value = 0.0943396226415094
$ws.Range("U2").Value = 0.00290275761973875
$ws.Range("W2").Value = 0.00145137880986938
$ws.Range("X2").Value = 0.06966618287373
$ws.Range("B3").Value = 0.870827285921626
$ws.Range("C3").Value = 0.971698113207547
$ws.Range("D3").Value = 0.0145137880986938
$ws.Range("E3").Value = 0.0188679245283019
$ws.Range("F3").Value = 0.000725689404934688
$ws.Range("G3").Value = 0.997822931785196
$ws.Range("H3").Value = 0.0587808417997097
$ws.Range("I3").Value = 0.0159651669085631
$ws.Range("J3").Value = 0.00290275761973875
$ws.Range("K3").Value = 0.153846153846154
$ws.Range("L3").Value = 0.963715529753266
$ws.Range("M3").Value = 0.00435413642960813
$ws.Range("N3").Value = 0.0341074020319303
$ws.Range("O3").Value = 0.991291727140784
$ws.Range("P3").Value = 0.00217706821480406
$ws.Range("Q3").Value = 0.933962264150943
$ws.Range("R3").Value = 0.0377358490566038
$ws.Range("S3").Value = 0.0217706821480406
$ws.Range("T3").Value = 0.00507982583454282
$ws.Range("U3").Value = 0.0188679245283019
$ws.Range("V3").Value = 0.979680696661829
$ws.Range("W3").Value = 0.979680696661829
$ws.Range("X3").Value = 0.910740203193033
$ws.Range("B4").Value = 0.032656023222061
$ws.Range("C4").Value = 0.0152394775036284
$ws.Range("D4").Value = 0.00217706821480406
$ws.Range("E4").Value = 0.111030478955007
$ws.Range("F4").Value = 0.997822931785196
$ws.Range("H4").Value = 0.0573294629898403
$ws.Range("I4").Value = 0.00725689404934688
$ws.Range("J4").Value = 0.00362844702467344
$ws.Range("K4").Value = 0.00507982583454282
$ws.Range("L4").Value = 0.00653120464441219
$ws.Range("M4").Value = 0.160377358490566
$ws.Range("N4").Value = 0.0377358490566038
$ws.Range("O4").Value = 0.00725689404934688
$ws.Range("P4").Value = 0.0058055152394775
$ws.Range("Q4").Value = 0.0195936139332366
$ws.Range("R4").Value = 0.529753265602322
$ws.Range("S4").Value = 0.021044992743106
$ws.Range("T4").Value = 0.887518142235123
$ws.Range("U4").Value = 0.976052249637155
$ws.Range("V4").Value = 0.00145137880986938
$ws.Range("W4").Value = 0.00217706821480406
$ws.Range("X4").Value = 0.00217706821480406
$ws.Range("B5").Value = 0.021044992743106
$ws.Range("C5").Value = 0.0108853410740203
$ws.Range("D5").Value = 0.983309143686502
$ws.Range("E5").Value = 0.0145137880986938
$ws.Range("G5").Value = 0.000725689404934688
$ws.Range("H5").Value = 0.0101596516690856
$ws.Range("I5").Value = 0.97677793904209
$ws.Range("J5").Value = 0.957184325108853
$ws.Range("K5").Value = 0.839622641509434
$ws.Range("M5").Value = 0.00362844702467344
$ws.Range("N5").Value = 0.920899854862119
$ws.Range("O5").Value = 0.000725689404934688
$ws.Range("P5").Value = 0.000725689404934688
$ws.Range("Q5").Value = 0.04644412191582
$ws.Range("R5").Value = 0.0101596516690856
$ws.Range("S5").Value = 0.934687953555878
$ws.Range("T5").Value = 0.0130624092888244
$ws.Range("U5").Value = 0.00145137880986938
$ws.Range("V5").Value = 0.0188679245283019
$ws.Range("W5").Value = 0.0166908563134978
$ws.Range("X5").Value = 0.0174165457184325
